$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes of a measure")

# Rename attribute labels in column A to wrap their unit suffixes in parentheses,
# e.g. "MEASURE.Exposure_s" -> "MEASURE.Exposure_(s)"
$ws.Range("A6").Value = "MEASURE.Exposure_(s)"
$ws.Range("A16").Value = "SPECTROMETER.Wavelength_(nm)"
$ws.Range("A17").Value = "SPECTROMETER.Confocal_Pinhole_Diameter_(AU)"
$ws.Range("A24").Value = "SPECTROMETER.Illumination_Power_(mW)"
$ws.Range("A30").Value = "SPECTROMETER.Scan_Amplitude_(GHz)"
$ws.Range("A32").Value = "SPECTROMETER.Scattering_Angle_(deg)"
$ws.Range("A33").Value = "SPECTROMETER.Spectral_Resolution_(MHz)"
$ws.Range("A34").Value = "SPECTROMETER.x-Mechanical_Resolution_(microm)"
$ws.Range("A35").Value = "SPECTROMETER.x-Optical_Resolution_(microm)"
$ws.Range("A36").Value = "SPECTROMETER.y-Mechanical_Resolution_(microm)"
$ws.Range("A37").Value = "SPECTROMETER.y-Optical_Resolution_(microm)"
$ws.Range("A38").Value = "SPECTROMETER.z-Mechanical_Resolution_(microm)"
$ws.Range("A39").Value = "SPECTROMETER.z-Optical_Resolution_(microm)"
$ws.Range("A10").Value = "MEASURE.Sampling_Step_Size_(dx,dy,dz)_(microm)"
$ws.Range("A11").Value = "MEASURE.Field_Of_View_(X,Y,Z)_(microm)"

# Update the saved selection on this sheet to A11 (matches the author's cursor
# position when the workbook was saved after the edit).
$ws.Range("A11").Select()
